# Update the "as_of_utc" timestamp column (AA) on the data sheets.
# Both "Главные" and "Линейные" carry the same stamp in AA2:AA26 and
# both need to move from 2025-11-28 03:05:41 -> 2025-11-28 07:06:03.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-28 07:06:03"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Cells.Item($row, 27).Value = $newTimestamp
    }
}
